# Update "Forecast Comparison" sheet with corrected forecast output:
#  - Insert a new "Week_Start_Date" column after "Week" (shifts all
#    subsequent columns one to the right, B..I -> C..J)
#  - Normalize "Week" labels from W01..W09 to W1..W9 (W10..W16 unchanged)
#  - Populate the new Week_Start_Date column with the week's start date

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before the current column B (ASIN), shifting the
# existing B:I columns to C:J.
$ws.Range("B1").EntireColumn.Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "Week_Start_Date"

# Week labels (column A) and their corresponding start dates (new column B).
$weeks = @(
    @{ Row = 2;  Week = "W1";  Date = "2025-01-05" },
    @{ Row = 3;  Week = "W2";  Date = "2025-01-12" },
    @{ Row = 4;  Week = "W3";  Date = "2025-01-19" },
    @{ Row = 5;  Week = "W4";  Date = "2025-01-26" },
    @{ Row = 6;  Week = "W5";  Date = "2025-02-02" },
    @{ Row = 7;  Week = "W6";  Date = "2025-02-09" },
    @{ Row = 8;  Week = "W7";  Date = "2025-02-16" },
    @{ Row = 9;  Week = "W8";  Date = "2025-02-23" },
    @{ Row = 10; Week = "W9";  Date = "2025-03-02" },
    @{ Row = 11; Week = "W10"; Date = "2025-03-09" },
    @{ Row = 12; Week = "W11"; Date = "2025-03-16" },
    @{ Row = 13; Week = "W12"; Date = "2025-03-23" },
    @{ Row = 14; Week = "W13"; Date = "2025-03-30" },
    @{ Row = 15; Week = "W14"; Date = "2025-04-06" },
    @{ Row = 16; Week = "W15"; Date = "2025-04-13" },
    @{ Row = 17; Week = "W16"; Date = "2025-04-20" }
)

foreach ($w in $weeks) {
    $ws.Cells.Item($w.Row, 1).Value = $w.Week

    # Force the date column to be stored as plain text (matching the
    # source data, which is an inline string rather than a date value).
    $dateCell = $ws.Cells.Item($w.Row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $w.Date
}
